# Commit: "Assigning Location functions to testers."
#
# This script:
#   1. Fills in a new "Tester" column (E) on the "Functions" sheet for the
#      "Location" class rows (rows 2-15), assigning each function to one of
#      three testers: Mel, Scott, or Nick.
#   2. Moves the active tab / selection from the "SW Units" sheet to the
#      "Functions" sheet (the author was last looking at the Functions tab,
#      with the cursor resting on C11 back on "SW Units").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "SW Units"
$ws2 = $wb.Worksheets.Item(2)   # "Functions"

# --- 1. Assign testers to the Location class functions (rows 2-15) -------
$testers = @{
    2  = "Mel"
    3  = "Mel"
    4  = "Mel"
    5  = "Mel"
    6  = "Scott"
    7  = "Scott"
    8  = "Nick"
    9  = "Scott"
    10 = "Nick"
    11 = "Nick"
    12 = "Mel"
    13 = "Scott"
    14 = "Nick"
    15 = "Scott"
}

foreach ($row in 2..15) {
    $ws2.Range("E$row").Value = $testers[$row]
}

# --- 2. Update the active sheet / selections ------------------------------
# "SW Units" keeps focus on C11 but is no longer the active tab.
$ws1.Range("C11").Select()

# "Functions" becomes the active tab.
$ws2.Activate()
$ws2.Range("A1").Select()
